$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows appended to the log (rows 306-322).
$colA = @('2026-02-11 18:25:33', '2026-02-11 15:17:49', '2026-02-11 12:24:19', '2026-02-11 15:22:13', '2026-02-11 14:00:57', '2026-02-11 13:42:44', '2026-02-11 11:42:19', '2026-02-11 17:34:56', '2026-02-11 09:56:00', '2026-02-11 14:59:28', '2026-02-11 11:21:14', '2026-02-11 12:06:11', '2026-02-11 13:47:47', '2026-02-11 14:40:19', '2026-02-11 14:49:33', '2026-02-11 17:31:02', '2026-02-11 18:21:36')
$colB = @('237676840777', '237677833877', '237678854978', '237679422591', '237650353920', '237651927448', '237653294562', '237678046498', '237679428698', '237679551262', '237680574202', '237681118330', '237674446293', '237679085953', '237681662761', '237682975726', '237683075075')
$colC = @('IVANS FANWOUM NOUPOUEH', 'MEDJOM TAGNE MICHELLE GUILENE LA NEGRESSE SARL', 'NSAMO NDJOUOHOU MICRANGE ETS MOBILE FINANCIAL SERVICES MFS', 'ETS LE CONTENT 42', 'MENIAPI HELENE EDOSSINE TOP MOBIL TELECOM', 'ODETTE MABAKOU EPOUSE KENNE', 'NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS', 'MFS SOCAVER', 'ETS LE CONTENT 29', 'LA NEGRESSE LTDLA CBOX R1 MEGNE JUDITH', 'TOUMEWO SAMUEL', 'SAHA NDESA JONAS LTDLA_POLAS_OTH_NDOGBONG SERIE', 'ARSENE TITCHO KWAKEP', 'FERNANDEZ NJOFANG TCHIYADJE', 'ROLCHILE DJAMEN KOUDJOU', 'LA NEGRESSE SARL NYOUNG JOSEPH CLOTAIRE', 'GLADYS LANG NGOINSEH')
$colD = @(193917, 131184, 92731, 3900, 612128, 172411, 847745, 6477, 42, 152434, 224117, 132896, 11136, 58169, 8697, 114072, 134739)

# Column B holds phone numbers; force text format so long numeric
# strings aren't coerced into numbers (matches existing columns A-C).
$ws.Range("B306:B322").NumberFormat = "@"

$startRow = 306
for ($i = 0; $i -lt $colA.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}
